$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.583.48'
$ws.Range('E2').Value = '  +1.82%  '
$ws.Range('D3').Value = '1.900.89'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').Value = '''244.53'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.20%  '
$ws.Range('D6').Value = '''0.634'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.05%  '
$ws.Range('E7').Value = '  +0.41%  '
$ws.Range('D8').Value = '''42.48'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('D10').Value = '''0.0706'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').Value = '''0.0994'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.98%  '
$ws.Range('D12').Value = '2.176.33'
$ws.Range('E12').Value = '  +3.09%  '
$ws.Range('D13').Value = '''12.50'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +9.90%  '
$ws.Range('D14').Value = '1.903.21'
$ws.Range('E14').Value = '  +3.29%  '
$ws.Range('E15').Value = '  +3.00%  '
$ws.Range('D16').Value = '''4.81'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.09%  '
$ws.Range('D17').Value = '35.554.05'
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('D18').Value = '''72.12'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.10%  '
$ws.Range('D19').Value = '0.0₃0809'
$ws.Range('E19').Value = '  +2.43%  '
$ws.Range('D20').Value = '''244.47'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('E21').Value = '  +2.32%  '
$ws.Range('D22').Value = '''4.91'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +3.12%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('D25').Value = '''171.10'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.43%  '
$ws.Range('D26').Value = '''2.10'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +27.39%  '
$ws.Range('D27').Value = '''8.45'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +8.26%  '
$ws.Range('D28').Value = '''17.96'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.90%  '
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('D30').Value = '''0.0566'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.65%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''4.09'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.61%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''0.947'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +27.53%  '
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('E34').Value = '  +5.44%  '
$ws.Range('E35').Value = '  +8.46%  '
$ws.Range('E36').Value = '  +4.52%  '
$ws.Range('E37').Value = '  +8.72%  '
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('D39').Value = '''0.0205'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.92%  '
$ws.Range('D40').Value = '''91.29'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.81%  '
$ws.Range('D41').Value = '1.357.94'
$ws.Range('E41').Value = '  +0.75%  '
$ws.Range('D42').Value = '''15.26'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.18%  '
$ws.Range('E43').Value = '  +13.23%  '
$ws.Range('D44').Value = '''13.08'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +32.09%  '
$ws.Range('D45').Value = '''2.35'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.66%  '
$ws.Range('D46').Value = '''46.94'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +38.81%  '
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('E48').Value = '  +0.42%  '
$ws.Range('D49').Value = '''6.67'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.50%  '
$ws.Range('D50').Value = '2.084.78'
$ws.Range('E50').Value = '  +2.73%  '
$ws.Range('E51').Value = '  +3.43%  '
